$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.017.06'
$ws.Range('E2').Value = '  -3.74%  '
$ws.Range('D3').Value = '2.974.88'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '557.12'
$ws.Range('E5').Value = '  -3.90%  '
$ws.Range('D6').Value = '134.02'
$ws.Range('E6').Value = '  +6.63%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +3.47%  '
$ws.Range('D9').Value = '2.963.48'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('E10').Value = '  -2.68%  '
$ws.Range('D11').Value = '4.86'
$ws.Range('E11').Value = '  -4.68%  '
$ws.Range('E12').Value = '  +2.50%  '
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').Value = '33.03'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '3.464.38'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = '6.83'
$ws.Range('E17').Value = '  +9.94%  '
$ws.Range('D18').Value = '2.963.16'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').Value = '57.868.95'
$ws.Range('E19').Value = '  -3.87%  '
$ws.Range('D20').Value = '419.54'
$ws.Range('E20').Value = '  -2.71%  '
$ws.Range('D21').Value = '13.22'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('E22').Value = '  +3.85%  '
$ws.Range('D23').Value = '6.99'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '13.04'
$ws.Range('E24').Value = '  +2.41%  '
$ws.Range('D25').Value = '79.60'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '0.997'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').Value = '7.60'
$ws.Range('E29').Value = '  +4.38%  '
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('D31').Value = '25.32'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '6.08'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').Value = '0.101'
$ws.Range('E33').Value = '  +7.74%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = '2.14'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '5.66'
$ws.Range('E35').Value = '  +1.45%  '
$ws.Range('D36').Value = '0.945'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('D37').Value = '0.0₃0694'
$ws.Range('E37').Value = '  +5.38%  '
$ws.Range('D38').Value = '48.60'
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('D39').Value = '8.52'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('D40').Value = '2.58'
$ws.Range('E40').Value = '  +4.98%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '382.69'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.109'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').Value = '2.683.91'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D46').Value = '0.243'
$ws.Range('E46').Value = '  +3.27%  '
$ws.Range('D47').Value = '122.67'
$ws.Range('E47').Value = '  +3.68%  '
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('D50').Value = '23.63'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('E51').Value = '  -0.22%  '
